$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Moksha")

# Column M (rows 2-131) holds example-comment text that belongs next to the
# free-text remark column S. Move the values (and their formatting) from M
# into S, leaving M blank.
$src = $ws.Range("M2:M131")
$dst = $ws.Range("S2:S131")
$src.Cut($dst)

# Leave the view the way the author left it: scrolled right so column H is
# the first visible column, with the now-empty source range selected.
$ws.Range("M2").Select()
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("M2:M131").Select()
